# Rename the "Kitimat Stikine" region label to the correctly hyphenated
# "Kitimat-Stikine" everywhere it appears as a Member value (column A),
# i.e. rows 52-59 of the msw_in_region_disposal sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("msw_in_region_disposal")

$oldName = "Kitimat Stikine"
$newName = "Kitimat-Stikine"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value -eq $oldName) {
        $cell.Value = $newName
    }
}

# Mirror the cursor/viewport position recorded in the saved file.
$ws.Range("A86").Select()
$excel.ActiveWindow.ScrollRow = 72
$excel.ActiveWindow.ScrollColumn = 1
